$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows for RM 232 and SC 92 (data points dropped from the dataset)
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# Apply corrected / imputed values for the error-calculation columns
$ws.Range("D3").Value = -14.2
$ws.Range("E3").Value = -5.7
$ws.Range("E4").Value = $null
$ws.Range("F4").Value = $null
$ws.Range("D5").Value = $null
$ws.Range("E5").Value = $null
$ws.Range("F5").Value = $null
$ws.Range("E9").Value = -6.8
$ws.Range("F9").Value = 17.26
$ws.Range("E10").Value = -6.1
$ws.Range("F10").Value = 16.43
$ws.Range("E11").Value = -7.9
$ws.Range("F11").Value = 17.65
$ws.Range("E12").Value = -5.3
$ws.Range("E15").Value = $null
$ws.Range("F15").Value = $null
$ws.Range("E17").Value = $null
$ws.Range("E18").Value = $null
$ws.Range("F18").Value = $null
$ws.Range("E20").Value = $null
$ws.Range("F20").Value = $null
$ws.Range("D21").Value = -14.3
$ws.Range("E21").Value = -8.699999999999999
$ws.Range("D23").Value = $null
$ws.Range("E23").Value = $null
$ws.Range("F23").Value = $null
$ws.Range("F25").Value = $null
$ws.Range("F27").Value = 17
$ws.Range("F28").Value = 17.44
$ws.Range("E31").Value = -8.1
$ws.Range("F31").Value = 17.18
$ws.Range("D32").Value = -14.7
$ws.Range("E32").Value = -6.4
$ws.Range("F32").Value = 17.39
